# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Thu Mar 28 03:06:50 UTC 2024 with GitHub Actions".
#
# Price (D) and Volume(1h) (E) cells are stored as literal text in the workbook
# (e.g. "3.475.19", "  -1.92%  "), so we force Text number format before writing
# into any cell whose new value could otherwise be auto-parsed by Excel as a
# number (e.g. "580.64" -> 580.64). Coin name (B) / link (C) cells are plain text
# and are written as-is; three rows are reordered in the ranking (27<->28 and the
# 35/36/37 rotation), so B/C/D/E are all rewritten for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='69.210.24'; E='  -1.95%  ' },
    @{ Row=3; D='3.482.09' },
    @{ Row=4; E='  -0.05%  ' },
    @{ Row=5; D='580.64'; E='  -0.61%  ' },
    @{ Row=6; D='182.14'; E='  -4.54%  ' },
    @{ Row=7; D='3.468.03'; E='  -4.08%  ' },
    @{ Row=8; D='0.607'; E='  -4.07%  ' },
    @{ Row=9; E='  +0.04%  ' },
    @{ Row=10; D='0.196'; E='  +6.32%  ' },
    @{ Row=11; D='0.639'; E='  -4.01%  ' },
    @{ Row=12; D='53.51'; E='  -4.83%  ' },
    @{ Row=13; D='0.0000303'; E='  -3.59%  ' },
    @{ Row=14; D='9.37'; E='  -3.94%  ' },
    @{ Row=15; D='4.020.81'; E='  -4.14%  ' },
    @{ Row=16; D='19.11'; E='  -4.71%  ' },
    @{ Row=17; D='69.076.95'; E='  -2.08%  ' },
    @{ Row=18; D='3.465.01'; E='  -4.09%  ' },
    @{ Row=19; D='12.20'; E='  -4.22%  ' },
    @{ Row=20; E='  -1.67%  ' },
    @{ Row=21; D='535.02'; E='  +8.90%  ' },
    @{ Row=22; E='  -4.79%  ' },
    @{ Row=23; D='18.33'; E='  -5.35%  ' },
    @{ Row=24; D='4.50'; E='  +2.16%  ' },
    @{ Row=25; D='4.85'; E='  -1.89%  ' },
    @{ Row=26; D='95.29'; E='  -1.67%  ' },
    @{ Row=27; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='2.94'; E='  -2.07%  ' },
    @{ Row=28; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='10.98'; E='  -1.26%  ' },
    @{ Row=29; D='9.03'; E='  -4.79%  ' },
    @{ Row=30; D='31.65'; E='  -2.48%  ' },
    @{ Row=31; D='7.17'; E='  -5.86%  ' },
    @{ Row=32; D='12.42'; E='  +0.96%  ' },
    @{ Row=33; D='63.54' },
    @{ Row=34; E='  -6.14%  ' },
    @{ Row=35; B='Bittensor'; C='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D='529.58'; E='  -8.92%  ' },
    @{ Row=36; B='Fetch.AI'; C='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D='3.08'; E='  +5.18%  ' },
    @{ Row=37; B='TheGraph'; C='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D='0.403'; E='  +0.31%  ' },
    @{ Row=38; E='  -0.17%  ' },
    @{ Row=39; D='37.64'; E='  -3.97%  ' },
    @{ Row=40; D='0.0₃0755'; E='  -8.13%  ' },
    @{ Row=41; E='  -3.56%  ' },
    @{ Row=42; E='  -2.53%  ' },
    @{ Row=43; D='3.328.27'; E='  +2.46%  ' },
    @{ Row=44; D='3.05'; E='  -6.23%  ' },
    @{ Row=45; D='3.47'; E='  +2.51%  ' },
    @{ Row=46; E='  -4.71%  ' },
    @{ Row=47; D='0.0435'; E='  -2.76%  ' },
    @{ Row=48; E='  -4.14%  ' },
    @{ Row=49; D='8.87'; E='  -8.81%  ' },
    @{ Row=50; D='0.999'; E='  +0.03%  ' },
    @{ Row=51; D='136.47'; E='  -1.71%  ' }
)

foreach ($u in $updates) {
    foreach ($col in @("B", "C", "D", "E")) {
        if (-not $u.ContainsKey($col)) { continue }
        $cell = $ws.Range("$col$($u.Row)")
        $value = $u[$col]
        if ($col -eq "D" -or $col -eq "E") {
            # Keep these columns as text even when the new value parses as a number
            # or looks like a date (e.g. "4.50", "0.999"), matching the inlineStr
            # storage already used throughout this sheet.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $value
    }
}
